$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before current row 11 (the maf_val row), shifting
# maf_val and clonal_threshold rows down by one.
$ws.Rows.Item(11).Insert()

# Row 10 changes from "missingness" to "locus_missingness" (value stays 0.4).
$ws.Cells.Item(10, 1).Value = "locus_missingness"

# New row 11: sample_missingness
# (Inserting the row already copied the B column style from row 10,
# matching the original workbook's style index for B11.)
$ws.Cells.Item(11, 1).Value = "sample_missingness"
$ws.Cells.Item(11, 2).Value = 0.2

# Description text updates (C11 then C10) to match new missingness split.
$ws.Cells.Item(11, 3).Value = "remove samples with missingness higher than this value"
$ws.Cells.Item(10, 3).Value = "remove loci with missingness higher than this value"

# Update the selected cell to reflect the new active cell.
$ws.Range("B11").Select()
